$d = $word.ActiveDocument

$replacements = @(
    @{old="333÷2="; new="637÷5="},
    @{old="946÷8="; new="113÷6="},
    @{old="288÷3="; new="686÷9="},
    @{old="205÷5="; new="533÷9="},
    @{old="490÷9="; new="628÷2="},
    @{old="313÷5="; new="796÷3="},
    @{old="846÷3="; new="118÷3="},
    @{old="127÷5="; new="145÷7="},
    @{old="913÷5="; new="492÷2="},
    @{old="124÷3="; new="300÷2="},
    @{old="623÷6="; new="689÷6="},
    @{old="571÷5="; new="285÷7="},
    @{old="642÷6="; new="586÷7="},
    @{old="779÷5="; new="763÷2="},
    @{old="506÷3="; new="113÷7="},
    @{old="647÷6="; new="757÷8="},
    @{old="261÷3="; new="631÷3="},
    @{old="477÷9="; new="450÷6="},
    @{old="462÷5="; new="144÷3="},
    @{old="464÷7="; new="779÷2="},
    @{old="280÷4="; new="201÷6="},
    @{old="168÷4="; new="297÷4="},
    @{old="177÷4="; new="676÷3="},
    @{old="868÷2="; new="342÷2="},
    @{old="449÷3="; new="732÷9="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
